$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Step 1: delete old rows 35 and 36 (Activation of Enemy Units / Placement of Enemy Units content)
# which collapses old rows 37,38,39 up to 35,36,37
$ws.Rows(35).Delete()
$ws.Rows(35).Delete()

# Step 2: insert a new blank row at position 38 to hold the new "Evening Debriefing" (e050) event
$ws.Rows(38).Insert()

# Step 3: write final content + row heights for rows 31-38

$ws.Range("A31").Value = "e030"
$b31 = @"
<Bold>e030 Advancing Fire Ammo Use</Bold> 
<InlineUIContainer><Button Content='r22.11' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Mark off 1D/2 (round down) HE rounds and .30 caliber MG ammo boxes regardless of whether the battle occurs or not. Mark off on the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
"@
$ws.Range("B31").Value = $b31
$ws.Rows(31).RowHeight = 105

$ws.Range("A32").Value = "e031"
$b32 = @"
<Bold>e031 Enemy Strength Roll Entering Battle Board</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table for enemy strength: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
"@
$ws.Range("B32").Value = $b32
$ws.Rows(32).RowHeight = 90

$ws.Range("A33").Value = "e032"
$b33 = @"
<Bold>e032 Battle Check</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to determine if combat occurs in this area: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
"@
$ws.Range("B33").Value = $b33
$ws.Rows(33).RowHeight = 90

$ws.Range("A34").Value = "e033"
$b34 = @"
<Bold>e033 No Combat</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
If converting territory to US Control, Victory points are added to the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
"@
$ws.Range("B34").Value = $b34
$ws.Rows(34).RowHeight = 75

$ws.Range("A35").Value = "e034"
$b35 = @"
<Bold>e034 Placing Advancing Fire Markers</Bold> 
<InlineUIContainer><Button Content='r4.61' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Place Advancing Fire Markers available to you per 
<InlineUIContainer><Button Content='r22.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
Click one of highlighted regions to place. 
<LineBreak/><LineBreak/>
You place up to six minus one marker for every three friendly tank losses (rounded up) . You may place more than one in a zone. 
<LineBreak/><LineBreak/>
The status bar on the bottom shows how many are remaining to place.
"@
$ws.Range("B35").Value = $b35
$ws.Rows(35).RowHeight = 150

$ws.Range("A36").Value = "e035"
$b36 = @"
<Bold>e035 Ambush Check</Bold> 
<InlineUIContainer><Button Content='r4.65' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D for possible Ambush: 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>&lt; 8
<LineBreak/><LineBreak/>
"@
$ws.Range("B36").Value = $b36
$ws.Rows(36).RowHeight = 90

$ws.Range("A37").Value = "e036"
$b37 = @"
<Bold>e036 Battle Board Empty</Bold> 
<InlineUIContainer><Button Content='r4.77' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Since the Battle Board is now empty of enemy units, the battle for this area is over. 
<LineBreak/><LineBreak/>
1.) Flip Resistance marker to US Controlled on Movement Board.
<LineBreak/>
2.) Victory points for control of the area added to the 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
3.) If daylight remains, return to Prepare for Battle per 
<InlineUIContainer><Button Content='r4.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
4.) No daylight, perform the Evening Debriefing per 
<InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
5.) Click image to continue.
<LineBreak/><LineBreak/>
          <InlineUIContainer><Image Name='Debrief' Height='225' Width='450'></Image></InlineUIContainer>
"@
$ws.Range("B37").Value = $b37
$ws.Rows(37).RowHeight = 285

$ws.Range("A38").Value = "e050"
$b38 = @"
<Bold>e050 Evening Debriefing</Bold> 
<LineBreak/><LineBreak/>
An evening debriefing is performed per <InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
Click image to continue.
<LineBreak/><LineBreak/>
                     <InlineUIContainer><Image Name='Sherman4' Height='168' Width='275'></Image></InlineUIContainer>
"@
$ws.Range("B38").Value = $b38
$ws.Rows(38).RowHeight = 90

# Step 4: update selection to match target view (B37 active cell)
$ws.Range("B37").Select()

Write-Output "done"
